$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly report data between row 2 and row 3
# Row 2 (Fecha, Volumen, Precio minimo, Precio promedio ponderado, Precio $/Kg)
$d2 = $ws.Range("D2").Value()
$j2 = $ws.Range("J2").Value()
$k2 = $ws.Range("K2").Value()
$m2 = $ws.Range("M2").Value()
$p2 = $ws.Range("P2").Value()

$d3 = $ws.Range("D3").Value()
$j3 = $ws.Range("J3").Value()
$k3 = $ws.Range("K3").Value()
$m3 = $ws.Range("M3").Value()
$p3 = $ws.Range("P3").Value()

$ws.Range("D2").Value = $d3
$ws.Range("J2").Value = $j3
$ws.Range("K2").Value = $k3
$ws.Range("M2").Value = $m3
$ws.Range("P2").Value = $p3

$ws.Range("D3").Value = $d2
$ws.Range("J3").Value = $j2
$ws.Range("K3").Value = $k2
$ws.Range("M3").Value = $m2
$ws.Range("P3").Value = $p2
